$d = $word.ActiveDocument

$replacements = @(
    @("406×4=1624", "567×2=1134"),
    @("354×4=1416", "284×7=1988"),
    @("773×7=5411", "281×2=562"),
    @("458×5=2290", "800×6=4800"),
    @("358×4=1432", "782×4=3128"),
    @("269×9=2421", "724×2=1448"),
    @("244×6=1464", "603×5=3015"),
    @("852×8=6816", "569×4=2276"),
    @("250×6=1500", "884×8=7072"),
    @("606×7=4242", "747×5=3735"),
    @("602×9=5418", "807×2=1614"),
    @("232×2=464", "167×3=501"),
    @("779×8=6232", "700×8=5600"),
    @("920×8=7360", "730×9=6570"),
    @("932×8=7456", "855×6=5130"),
    @("418×6=2508", "661×3=1983"),
    @("952×6=5712", "835×8=6680"),
    @("646×7=4522", "109×7=763"),
    @("124×5=620", "885×2=1770"),
    @("596×7=4172", "871×2=1742"),
    @("173×5=865", "758×8=6064"),
    @("251×7=1757", "310×8=2480"),
    @("997×4=3988", "419×9=3771"),
    @("647×3=1941", "948×6=5688"),
    @("797×5=3985", "268×7=1876")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
